$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndividualBiometrics")

# Add two new header columns: H1 = "Protein", I1 = "Ontogeny"
$ws.Range("H1").Value = "Protein"
$ws.Range("I1").Value = "Ontogeny"

# Update the selection to I2
$ws.Range("I2").Select()
